# Apply question/answer text corrections to the Erste-Hilfe quiz sheet,
# matching commit "Questions corrected, correctAnswer highlighted".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Herzinfarkt): fix "Herzbereich" -> "Herzbereich/Brustkorb" and typo "Obernauch" -> "Oberbauch"
$ws.Range("E6").Value = "Starke Schmerzen in der Brust, Engegefühl im Herzbereich/Brustkorb, Übelkeit, Erbrechen, Atemnot und Schmerzen im Oberbauch, fahles Gesicht"

# Row 7 (stabile Seitenlage): clarify correct answer
$ws.Range("D7").Value = "Um die Atemwege frei zu machen und damit Blut und Erbrochenes abfließen können und die bewusstlose Person nicht daran erstickt"

# Row 8 (Herzstillstand): replace wrong-answer text and clarify correct answer
$ws.Range("B8").Value = "Da kann man nichts machen"
$ws.Range("E8").Value = "Notruf absetzen, Hilfe holen, Defi holen lassen, 30x Herzdruckmassage, 2x Mund-zu-Mund-Beatmung im Wechsel"

# Row 21 (Rautek-Rettungsgriff): add clarification about thumb position
$ws.Range("D21").Value = "Unter den Achseln der Person durchgreifen, den Unterarm von oben fassen (Daumen oben lassen!), die Person anheben und wegschleifen"

# Row 22 (Nasenbluten): fix typo "Taschntüchern" -> "Taschentüchern"
$ws.Range("E22").Value = "Nasenlöcher mit vielen Taschentüchern zustopfen"

# Row 26 (Fleisch verschluckt): "Ein Erwachsener" -> "Eine Person"
$ws.Range("A26").Value = "Eine Person hat ein Stück Fleisch verschluckt und kann weder atmen noch sprechen. Was sollten Sie als Ersthelfender unverzüglich tun?"

# Move the visible selection to G26, matching the saved selection in the workbook.
$ws.Range("G26").Select()
